# Apply updated cryptocurrency price/volume data as of Sat Jul 29 20:19:43 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.362.68"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.881.82"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7142"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08068"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3140"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08356"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.859.78"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7204"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.284"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008402"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "29.354.61"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "2.120.26"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1591"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.083"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.429"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.353"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.214"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05381"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.954"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7532"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.180"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("D39").Value = "1.280.61"
$ws.Range("E39").Value = "  +9.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8900"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000130"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.78%  "
$ws.Range("D47").Value = "2.018.85"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.805"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5208"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.477"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4371"
$ws.Range("D51").Style = "Normal"
